$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted right before the existing
# row 139 (the data block starts at row 2; row 139 is currently the record
# dated 2021-03-24 / D=44279). Inserting the row shifts every following
# record down by one (old row 139 -> new row 140, ..., old row 222 -> new
# row 223), which is exactly what the target diff shows.
$ws.Rows.Item(139).Insert()

# Fill in the new record. Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant for
# every row in this sub-table, so we just replicate them; D,J,K,L,M,P hold
# the new week's data.
$ws.Range("A139").Value = 3
$ws.Range("B139").Value = "Femacal de La Calera"
$ws.Range("C139").Value = "Coquimbo"
$ws.Range("D139").Value = 44830
$ws.Range("E139").Value = 5
$ws.Range("F139").Value = 100112010
$ws.Range("G139").Value = "Achicoria"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 60
$ws.Range("K139").Value = 6000
$ws.Range("L139").Value = 6000
$ws.Range("M139").Value = 6000
$ws.Range("N139").Value = "$/caja 16 unidades"
$ws.Range("O139").Value = "Provincia de Quillota"
$ws.Range("P139").Value = 375
$ws.Range("Q139").Value = 16
$ws.Range("R139").Value = "Hortaliza"
